$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 315.125
$ws.Range("I2").Value = 115.5
$ws.Range("J2").Value = 914
$ws.Range("K2").Value = 115.5
$ws.Range("L2").Value = 914
$ws.Range("M2").Value = -2.5
$ws.Range("N2").Value = -1140

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 142
$ws.Range("I9").Value = 172
$ws.Range("J9").Value = 124
$ws.Range("K9").Value = 172
$ws.Range("L9").Value = 124
$ws.Range("M9").Value = -3
$ws.Range("N9").Value = -462

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 219
$ws.Range("I33").Value = 171
$ws.Range("K33").Value = 171
$ws.Range("M33").Value = 58

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 10900.1
$ws.Range("J51").Value = 9888.888999999999
$ws.Range("L51").Value = 9888.888999999999
$ws.Range("N51").Value = -10856.889

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H127").Value = 3656.6667
$ws.Range("I127").Value = 3656.6667
$ws.Range("J127").Value = 0
$ws.Range("K127").Value = 10970.0001
$ws.Range("L127").Value = 0
$ws.Range("M127").Value = -6010.000100000001
$ws.Range("N127").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2189.1404
$ws.Range("I138").Value = 1509.3334
$ws.Range("J138").Value = 2431.9285
$ws.Range("K138").Value = 4528.0002
$ws.Range("L138").Value = 7295.7855
$ws.Range("M138").Value = 611.9997999999996
$ws.Range("N138").Value = -17575.7855

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 7399.4
$ws.Range("I61").Value = 7999.25
$ws.Range("J61").Value = 5000
$ws.Range("K61").Value = 7999.25
$ws.Range("L61").Value = 5000
$ws.Range("M61").Value = -7787.25
$ws.Range("N61").Value = -5424

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2050.8
$ws.Range("I122").Value = 2054.8572
$ws.Range("K122").Value = 6164.571599999999
$ws.Range("M122").Value = -3714.571599999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 7399.4
$ws.Range("I136").Value = 7999.25
$ws.Range("J136").Value = 5000
$ws.Range("K136").Value = 23997.75
$ws.Range("L136").Value = 15000
$ws.Range("M136").Value = -21447.75
$ws.Range("N136").Value = -20100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1897.7142
$ws.Range("I20").Value = 2374.6667
$ws.Range("K20").Value = 2374.6667
$ws.Range("M20").Value = -2127.6667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 235.75
$ws.Range("I22").Value = 166.33333
$ws.Range("K22").Value = 166.33333
$ws.Range("M22").Value = 6.666670000000011

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2004.75
$ws.Range("I94").Value = 1576.8572
$ws.Range("K94").Value = 1576.8572
$ws.Range("M94").Value = -1125.8572

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1251.5238
$ws.Range("I134").Value = 1014.1053
$ws.Range("K134").Value = 3042.3159
$ws.Range("M134").Value = -507.3159000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 782.2
$ws.Range("I7").Value = 628
$ws.Range("J7").Value = 848.2857
$ws.Range("K7").Value = 628
$ws.Range("L7").Value = 848.2857
$ws.Range("M7").Value = -515
$ws.Range("N7").Value = -1074.2857

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3700.2
$ws.Range("I31").Value = 3079.7778
$ws.Range("J31").Value = 4630.8335
$ws.Range("K31").Value = 3079.7778
$ws.Range("L31").Value = 4630.8335
$ws.Range("M31").Value = -2784.7778
$ws.Range("N31").Value = -5220.8335

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3700.2
$ws.Range("I34").Value = 3079.7778
$ws.Range("J34").Value = 4630.8335
$ws.Range("K34").Value = 3079.7778
$ws.Range("L34").Value = 4630.8335
$ws.Range("M34").Value = -2877.7778
$ws.Range("N34").Value = -5034.8335

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 34897.5
$ws.Range("J51").Value = 34897.5
$ws.Range("L51").Value = 34897.5
$ws.Range("N51").Value = -36369.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H59").Value = 40882.668
$ws.Range("J59").Value = 41462.91
$ws.Range("L59").Value = 41462.91
$ws.Range("N59").Value = -43752.91

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H61").Value = 34897.5
$ws.Range("J61").Value = 34897.5
$ws.Range("L61").Value = 34897.5
$ws.Range("N61").Value = -35593.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 2253922.5
$ws.Range("I4").Value = 3755933.2
$ws.Range("K4").Value = 11267799.6
$ws.Range("M4").Value = -11267687.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 122.52381
$ws.Range("I12").Value = 120.2
$ws.Range("K12").Value = 360.6
$ws.Range("M12").Value = -187.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 5000335
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 5000335
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 15001005
$ws.Range("N46").Value = -15001187
$ws.Range("M46").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 17349.8
$ws.Range("I70").Value = 13699.2
$ws.Range("K70").Value = 13699.2
$ws.Range("M70").Value = -13429.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 17349.8
$ws.Range("I73").Value = 13699.2
$ws.Range("K73").Value = 13699.2
$ws.Range("M73").Value = -12763.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 902
$ws.Range("I102").Value = 804
$ws.Range("J102").Value = 1000
$ws.Range("K102").Value = 804
$ws.Range("L102").Value = 1000
$ws.Range("M102").Value = 818
$ws.Range("N102").Value = -4244

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 118562
$ws.Range("I122").Value = 6343.1665
$ws.Range("J122").Value = 342999.66
$ws.Range("K122").Value = 19029.4995
$ws.Range("L122").Value = 1028998.98
$ws.Range("N122").Value = -1033898.98
$ws.Range("M122").Value = -16579.4995

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 5309.857
$ws.Range("I126").Value = 4775
$ws.Range("J126").Value = 6023
$ws.Range("K126").Value = 14325
$ws.Range("L126").Value = 18069
$ws.Range("M126").Value = -11855
$ws.Range("N126").Value = -23009

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 78358.84
$ws.Range("I20").Value = 112833.89
$ws.Range("J20").Value = 790
$ws.Range("K20").Value = 112833.89
$ws.Range("L20").Value = 790
$ws.Range("M20").Value = -112607.89
$ws.Range("N20").Value = -1242

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 3667.6667
$ws.Range("J68").Value = 4501.5
$ws.Range("L68").Value = 4501.5
$ws.Range("N68").Value = -5999.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 3667.6667
$ws.Range("J71").Value = 4501.5
$ws.Range("L71").Value = 22507.5
$ws.Range("N71").Value = -29995.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 2526.8572
$ws.Range("I100").Value = 2448
$ws.Range("J100").Value = 3000
$ws.Range("K100").Value = 2448
$ws.Range("L100").Value = 3000
$ws.Range("M100").Value = -1907
$ws.Range("N100").Value = -4082

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 252000.5
$ws.Range("I132").Value = 252000.5
$ws.Range("K132").Value = 756001.5
$ws.Range("M132").Value = -753471.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 17618.25
$ws.Range("J45").Value = 17618.25
$ws.Range("L45").Value = 17618.25
$ws.Range("N45").Value = -18600.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("I113").Value = 800
$ws.Range("K113").Value = 2400
$ws.Range("M113").Value = -230

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 5292.8335
$ws.Range("I132").Value = 5219.4287
$ws.Range("J132").Value = 5549.75
$ws.Range("K132").Value = 15658.2861
$ws.Range("L132").Value = 16649.25
$ws.Range("M132").Value = -13128.2861
$ws.Range("N132").Value = -21709.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2562.0908
$ws.Range("I136").Value = 2506.158
$ws.Range("K136").Value = 7518.474
$ws.Range("M136").Value = -4968.474
